# Weekly refresh: insert a new data row at the top of the Plátano price
# series (row 251) pushing all later rows down by one, and carry the
# constant "identity" columns (Mercado/Región/Tipo/Producto/.../Unidad/
# Origen/Kg-unidad) across from the row that used to occupy row 251.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 251:332 down to 252:333, creating a blank row 251.
$ws.Rows(251).Insert()

# Populate the new row 251 with this week's record.
$ws.Range("A251").Value = 4
$ws.Range("B251").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C251").Value = "Los Lagos"
$ws.Range("D251").Value = 44559
$ws.Range("E251").Value = 10
$ws.Range("F251").Value = "Fruta"
$ws.Range("G251").Value = 100108
$ws.Range("H251").Value = "Tropicales y subtropicales"
$ws.Range("I251").Value = 100108006
$ws.Range("J251").Value = "Plátano"
$ws.Range("K251").Value = "Sin especificar"
$ws.Range("L251").Value = "Primera Pintón"
$ws.Range("M251").Value = 400
$ws.Range("N251").Value = 18000
$ws.Range("O251").Value = 19000
$ws.Range("P251").Value = 18500
$ws.Range("Q251").Value = "`$/caja 20 kilos"
$ws.Range("R251").Value = "Ecuador"
$ws.Range("S251").Value = 925
$ws.Range("T251").Value = 20
